$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 360
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
